$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (dSF) updates per diff
$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -15
$ws.Range("F5").Value = -4
$ws.Range("F7").Value = -6
$ws.Range("F9").Value = -3
$ws.Range("F11").Value = -1
$ws.Range("F14").Value = -8
$ws.Range("F17").Value = 4
